# EIA Table 1.2.B monthly update: October 2016 -> November 2016
# 1) Title + "Rolling 12 Months Ending in ..." labels
# 2) Insert a new "November" data row in the "Year 2016" block (IPP section)
# 3) Shift "Year to Date" and "Rolling 12 Months Ending in ..." blocks down by one row
# 4) Refresh the Year-to-Date and Rolling-12-months totals with the new figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the report title and the "Rolling 12 Months..." caption text.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value2 = "Table 1.2.B Net Generation by Energy Source:  Independent Power Producers, 2006-November 2016"

# ---------------------------------------------------------------------------
# 2. Shift rows 53-60 down to 54-61 (bottom-up so sources aren't clobbered),
#    carrying formatting (incl. row height for the footnote row) and values.
# ---------------------------------------------------------------------------
for ($src = 60; $src -ge 53; $src--) {
    $dst = $src + 1
    $srcRange = $ws.Range("A" + $src + ":M" + $src)
    $dstRange = $ws.Range("A" + $dst + ":M" + $dst)
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats
    $dstRange.Value2 = $srcRange.Value2
}
$ws.Application.CutCopyMode = 0

# Re-establish the merged header/footer rows at their new location.
$ws.Range("A53:M53").UnMerge()
$ws.Range("A57:M57").UnMerge()
$ws.Range("A60:M60").UnMerge()
$ws.Range("A54:M54").Merge()
$ws.Range("A58:M58").Merge()
$ws.Range("A61:M61").Merge()

# ---------------------------------------------------------------------------
# 3. Populate the new "November" row (row 53) with the same look as the
#    other month rows (copy format from the October row directly above).
# ---------------------------------------------------------------------------
$ws.Range("A52:M52").Copy()
$ws.Range("A53:M53").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = 0

$ws.Range("A53").Value2 = "November"
$ws.Range("B53").Value2 = 21420
$ws.Range("C53").Value2 = 338
$ws.Range("D53").Value2 = 116
$ws.Range("E53").Value2 = 42440
$ws.Range("F53").Value2 = 338
$ws.Range("G53").Value2 = 32097
$ws.Range("H53").Value2 = 1003
$ws.Range("I53").Value2 = 2458
$ws.Range("J53").Value2 = 19667
$ws.Range("K53").Value2 = -85
$ws.Range("L53").Value2 = 576
$ws.Range("M53").Value2 = 120368

# ---------------------------------------------------------------------------
# 4. Refresh "Year to Date" totals (rows 55-57: 2014, 2015, 2016).
# ---------------------------------------------------------------------------
function Set-RowValues($row, $values) {
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $row).Value2 = $values[$i]
    }
}

Set-RowValues 55 @(2014, 365427, 6561, 1271, 489106, 2898, 342360, 18074, 15147, 180622, -959, 6109, 1426616)
Set-RowValues 56 @(2015, 322545, 5987, 1533, 571163, 3214, 346861, 16275, 21514, 182104, -916, 6231, 1476513)
Set-RowValues 57 @(2016, 276328, 3060, 1276, 584558, 3598, 346505, 16381, 31237, 211850, -961, 6494, 1480326)

# ---------------------------------------------------------------------------
# 5. Update the "Rolling 12 Months Ending in ..." caption (row 58) and its
#    totals (rows 59-60: 2015, 2016).
# ---------------------------------------------------------------------------
$ws.Range("A58").Value2 = "Rolling 12 Months Ending in November"

Set-RowValues 59 @(2015, 352819, 6215, 1672, 613815, 3563, 381796, 18062, 22453, 198206, -986, 6812, 1604426)
Set-RowValues 60 @(2016, 296392, 3313, 1344, 633234, 3900, 380142, 18102, 32685, 232604, -1032, 7101, 1607784)

# ---------------------------------------------------------------------------
# 6. Row height is row-level (not cell-level), so none of the PasteSpecial /
#    value writes above touch it directly - fix it up last. Every shifted
#    row is default height except the footnote, which moved from row 60 to
#    row 61 and keeps its tall custom height.
# ---------------------------------------------------------------------------
$ws.Rows(60).RowHeight = 15
$ws.Rows(60).AutoFit()
$ws.Rows(61).RowHeight = 234
